$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 884, shifting existing rows (884-997) down to (885-998)
$ws.Rows.Item(884).EntireRow.Insert()

# Populate the newly inserted row 884 with the new data record
$ws.Range("A884").Value = 10
$ws.Range("B884").Value = "Vega Modelo de Temuco"
$ws.Range("C884").Value = "La Araucanía"
$ws.Range("D884").Value = 44946
$ws.Range("E884").Value = 9
$ws.Range("F884").Value = 100114001
$ws.Range("G884").Value = "Papa"
$ws.Range("H884").Value = "Patagonia"
$ws.Range("I884").Value = "1a (guarda)"
$ws.Range("J884").Value = 650
$ws.Range("K884").Value = 11000
$ws.Range("L884").Value = 11000
$ws.Range("M884").Value = 11000
$ws.Range("N884").Value = "`$/saco 25 kilos"
$ws.Range("O884").Value = "Provincia de Cautín"
$ws.Range("P884").Value = 440
$ws.Range("Q884").Value = 25
$ws.Range("R884").Value = "Hortaliza"
